# Matlab/Messreihen.xlsx - "programm anpassungen in Matlab, neue Messreihen"
#
# Adds two new measurement series (M18/M19) right after the existing M01..M17
# block, and their corresponding description rows (D18/D19) at the end of the
# D06..D17 block. Rows 35/36 were blank placeholder rows before this edit and
# become the new M18/M19 rows; a single blank separator row is kept before the
# "D" series (which is why that whole block shifts down by exactly one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "M" measurement rows (rows 35 and 36 were empty before) ----------
$ws.Cells.Item(35, 1).Value = "M18"
$ws.Cells.Item(35, 2).Value = "18K"
$ws.Cells.Item(35, 3).Value = 570
$ws.Cells.Item(35, 4).Value = 500
$ws.Cells.Item(35, 5).Value = -5.5
$ws.Cells.Item(35, 6).Value = "Schaltnetzteil"
$ws.Cells.Item(35, 7).Value = "3V LM317 Spannungs teiler von U nicht auf 3V ausgelegt"

$ws.Cells.Item(36, 1).Value = "M19"
$ws.Cells.Item(36, 2).Value = "2,7K"
$ws.Cells.Item(36, 3).Value = 100
$ws.Cells.Item(36, 4).Value = 500
$ws.Cells.Item(36, 5).Value = -5.5
$ws.Cells.Item(36, 6).Value = "Schaltnetzteil"
$ws.Cells.Item(36, 7).Value = "3V LM317 U jetzt wieder voll ausgesteuert"

# --- Keep a blank separator row before the "D" series ----------------------
# Previously the blank separator spanned rows 35-36 (before row 37 = D06).
# Now that 35/36 hold data, insert a fresh blank row so the separator still
# exists right before the D-series (D06 moves from row 37 to row 38, and so
# on through D17 moving from 48 to 49).
$ws.Rows.Item(37).Insert()

# --- New "D" description rows appended after D17 (now at row 49) ----------
$ws.Cells.Item(50, 1).Value = "D18"
$ws.Cells.Item(50, 2).Value = "18K"
$ws.Cells.Item(50, 3).Value = 570
$ws.Cells.Item(50, 4).Value = 500
$ws.Cells.Item(50, 5).Value = -5.5
$ws.Cells.Item(50, 6).Value = "Schaltnetzteil"
$ws.Cells.Item(50, 7).Value = "3V LM317 Spannungs teiler von U nicht auf 3V ausgelegt"

$ws.Cells.Item(51, 1).Value = "D19"
$ws.Cells.Item(51, 2).Value = "2,7K"
$ws.Cells.Item(51, 3).Value = 100
$ws.Cells.Item(51, 4).Value = 500
$ws.Cells.Item(51, 5).Value = -5.5
$ws.Cells.Item(51, 6).Value = "Schaltnetzteil"
$ws.Cells.Item(51, 7).Value = "3V LM317 U jetzt wieder voll ausgesteuert"

# --- Update the view state to match where the author ended up editing -----
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G56").Select()
